# Apply the commit: "Update gh-pages to output generated at 456a3b4"
# Refreshes event data across the four sheets: 展览, 演出, 本地生活, 全部类型
$wb = $excel.ActiveWorkbook

# ---------------- Sheet: 展览 (Exhibitions) ----------------
$ws1 = $wb.Worksheets.Item("展览")
$data_ws1 = New-Object 'object[,]' 20,8
$data_ws1[0,0] = '2024.03.02'
$data_ws1[0,1] = '广州·明日方舟ONLY'
$data_ws1[0,2] = '清河东路288号 科尔海悦酒店'
$data_ws1[0,3] = '2024.03.02 10:00-03.02 17:00'
$data_ws1[0,4] = 410
$data_ws1[0,5] = 60
$data_ws1[0,6] = 'https://show.bilibili.com/platform/detail.html?id=80271'
$data_ws1[0,7] = '//i0.hdslb.com/bfs/openplatform/202312/O9z5j7RB1703733578857.jpeg'
$data_ws1[1,0] = '2024.03.09'
$data_ws1[1,1] = '广州·VOCALOID术力口only'
$data_ws1[1,2] = '黄边三横路一街1号 设计殿堂'
$data_ws1[1,3] = '2024.03.09 09:30-03.09 18:00'
$data_ws1[1,4] = 553
$data_ws1[1,5] = 68.8
$data_ws1[1,6] = 'https://show.bilibili.com/platform/detail.html?id=81398'
$data_ws1[1,7] = '//i2.hdslb.com/bfs/openplatform/202401/XpsHJTsC1706160000879.png'
$data_ws1[2,0] = '2024.03.09'
$data_ws1[2,1] = '广州·排球少年.only'
$data_ws1[2,2] = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$data_ws1[2,3] = '2024.03.09 10:00-03.09 17:00'
$data_ws1[2,4] = 868
$data_ws1[2,5] = 60
$data_ws1[2,6] = 'https://show.bilibili.com/platform/detail.html?id=80453'
$data_ws1[2,7] = '//i1.hdslb.com/bfs/openplatform/202401/JYZdnQHl1704341212206.jpeg'
$data_ws1[3,0] = '2024.03.09'
$data_ws1[3,1] = '广州·樱漫动漫嘉年华9.0'
$data_ws1[3,2] = '奥体南路12号 优托邦(奥体旗舰店)'
$data_ws1[3,3] = '2024.03.09 10:00-03.09 17:00'
$data_ws1[3,4] = 595
$data_ws1[3,5] = 58
$data_ws1[3,6] = 'https://show.bilibili.com/platform/detail.html?id=81785'
$data_ws1[3,7] = '//i0.hdslb.com/bfs/openplatform/202402/gBrd7lLX1707203945413.jpeg'
$data_ws1[4,0] = '2024.03.10'
$data_ws1[4,1] = '广州·进击的巨人only'
$data_ws1[4,2] = '机场路1399号广州百信广场二期 李宁运动中心'
$data_ws1[4,3] = '2024.03.10 10:00-03.10 17:00'
$data_ws1[4,4] = 783
$data_ws1[4,5] = 60
$data_ws1[4,6] = 'https://show.bilibili.com/platform/detail.html?id=80454'
$data_ws1[4,7] = '//i0.hdslb.com/bfs/openplatform/202401/m8QkMfFl1704347300282.jpeg'
$data_ws1[5,0] = '2024.03.10'
$data_ws1[5,1] = '广州·马娘only2024'
$data_ws1[5,2] = '黄边三横路一街1号 设计殿堂'
$data_ws1[5,3] = '2024.03.10 09:30-03.10 18:00'
$data_ws1[5,4] = 354
$data_ws1[5,5] = 68.8
$data_ws1[5,6] = 'https://show.bilibili.com/platform/detail.html?id=81632'
$data_ws1[5,7] = '//i0.hdslb.com/bfs/openplatform/202402/svWCXpKm1706776489024.png'
$data_ws1[6,0] = '2024.03.16'
$data_ws1[6,1] = '广州·SISP动漫游戏嘉年华'
$data_ws1[6,2] = '西湾路150号 悦汇城'
$data_ws1[6,3] = '2024.03.16 13:00-03.17 19:00'
$data_ws1[6,4] = 541
$data_ws1[6,5] = 48
$data_ws1[6,6] = 'https://show.bilibili.com/platform/detail.html?id=80624'
$data_ws1[6,7] = '//i2.hdslb.com/bfs/openplatform/202401/Z4Q6Fv8B1704770554777.jpeg'
$data_ws1[7,0] = '2024.03.16'
$data_ws1[7,1] = '广州·YU 7th动漫嘉年华'
$data_ws1[7,2] = '珠江西路8号 高德置地夏广场'
$data_ws1[7,3] = '2024.03.16 10:00-03.17 17:00'
$data_ws1[7,4] = 104
$data_ws1[7,5] = 55
$data_ws1[7,6] = 'https://show.bilibili.com/platform/detail.html?id=81627'
$data_ws1[7,7] = '//i1.hdslb.com/bfs/openplatform/202402/lVqoZMVQ1706775042937.jpeg'
$data_ws1[8,0] = '2024.03.16'
$data_ws1[8,1] = '广州·代号鸢only2.0'
$data_ws1[8,2] = '清河东路288号 科尔海悦酒店'
$data_ws1[8,3] = '2024.03.16 10:00-03.16 21:00'
$data_ws1[8,4] = 1109
$data_ws1[8,5] = 39
$data_ws1[8,6] = 'https://show.bilibili.com/platform/detail.html?id=79828'
$data_ws1[8,7] = '//i0.hdslb.com/bfs/openplatform/202312/RVUVc8oy1702549585918.jpeg'
$data_ws1[9,0] = '2024.03.16'
$data_ws1[9,1] = '广州·原神X星穹铁道X绝区零ONLY'
$data_ws1[9,2] = '洛浦街夏滘西环路1号(厦滘地铁站A口步行290米) 厦喾岭南电商园会展中心'
$data_ws1[9,3] = '2024.03.16 10:00-03.16 17:00'
$data_ws1[9,4] = 558
$data_ws1[9,5] = 60
$data_ws1[9,6] = 'https://show.bilibili.com/platform/detail.html?id=80715'
$data_ws1[9,7] = '//i0.hdslb.com/bfs/openplatform/202401/Lt6ZYvA41704878219924.jpeg'
$data_ws1[10,0] = '2024.03.23'
$data_ws1[10,1] = '广州·BanG Dream ONLY'
$data_ws1[10,2] = '西环路1号 广州岭南会展中心'
$data_ws1[10,3] = '2024.03.23 10:00-03.23 17:00'
$data_ws1[10,4] = 322
$data_ws1[10,5] = 65
$data_ws1[10,6] = 'https://show.bilibili.com/platform/detail.html?id=81754'
$data_ws1[10,7] = '//i1.hdslb.com/bfs/openplatform/202402/3HJiKSeD1707104926306.jpeg'
$data_ws1[11,0] = '2024.03.23'
$data_ws1[11,1] = '广州·排球少年ONLY'
$data_ws1[11,2] = '机场路1399号广州百信广场二期 李宁运动中心'
$data_ws1[11,3] = '2024.03.23 10:00-03.23 17:00'
$data_ws1[11,4] = 445
$data_ws1[11,5] = 60
$data_ws1[11,6] = 'https://show.bilibili.com/platform/detail.html?id=80716'
$data_ws1[11,7] = '//i0.hdslb.com/bfs/openplatform/202401/IFLvYmxx1704879325152.jpeg'
$data_ws1[12,0] = '2024.03.24'
$data_ws1[12,1] = '广州·妖都恋与制作人ONLY3.0'
$data_ws1[12,2] = '迎宾大道123号 赛仑吉地大酒店'
$data_ws1[12,3] = '2024.03.24 10:00-03.24 17:30'
$data_ws1[12,4] = 135
$data_ws1[12,5] = 68
$data_ws1[12,6] = 'https://show.bilibili.com/platform/detail.html?id=81715'
$data_ws1[12,7] = '//i1.hdslb.com/bfs/openplatform/202402/SOOWI9wL1708675967102.jpeg'
$data_ws1[13,0] = '2024.04.06'
$data_ws1[13,1] = '广州·运动番only'
$data_ws1[13,2] = '机场路1399号广州百信广场二期 李宁运动中心'
$data_ws1[13,3] = '2024.04.06 10:00-04.06 17:00'
$data_ws1[13,4] = 287
$data_ws1[13,5] = 55
$data_ws1[13,6] = 'https://show.bilibili.com/platform/detail.html?id=81454'
$data_ws1[13,7] = '//i2.hdslb.com/bfs/openplatform/202401/TBZfwnB41706255329549.jpeg'
$data_ws1[14,0] = '2024.04.13'
$data_ws1[14,1] = '广州·Veni Vidi Vici动漫游戏嘉年华'
$data_ws1[14,2] = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$data_ws1[14,3] = '2024.04.13 10:00-04.13 17:00'
$data_ws1[14,4] = 41
$data_ws1[14,5] = 58
$data_ws1[14,6] = 'https://show.bilibili.com/platform/detail.html?id=81575'
$data_ws1[14,7] = '//i1.hdslb.com/bfs/openplatform/202401/7ir7DZHt1706697841803.jpeg'
$data_ws1[15,0] = '2024.04.13'
$data_ws1[15,1] = '广州·潮娃展WWS'
$data_ws1[15,2] = '西环路1号 广州岭南会展中心'
$data_ws1[15,3] = '2024.04.13 10:00-04.13 17:00'
$data_ws1[15,4] = 64
$data_ws1[15,5] = 48
$data_ws1[15,6] = 'https://show.bilibili.com/platform/detail.html?id=81745'
$data_ws1[15,7] = '//i2.hdslb.com/bfs/openplatform/202402/1SWNaBsA1707100228293.jpeg'
$data_ws1[16,0] = '2024.04.20'
$data_ws1[16,1] = '广州·Arknights Only·夜航星（明日方舟Only)'
$data_ws1[16,2] = '同泰路颐和山庄 颐和大酒店'
$data_ws1[16,3] = '2024.04.20 10:00-04.20 17:00'
$data_ws1[16,4] = 516
$data_ws1[16,5] = 69
$data_ws1[16,6] = 'https://show.bilibili.com/platform/detail.html?id=80282'
$data_ws1[16,7] = '//i2.hdslb.com/bfs/openplatform/202312/gaEHIE1F1703745559785.jpeg'
$data_ws1[17,0] = '2024.05.05'
$data_ws1[17,1] = '广州·第八届萌物语动漫嘉年华'
$data_ws1[17,2] = '洛浦街厦滘西环路1号 岭南会展中心'
$data_ws1[17,3] = '2024.05.05 10:00-05.05 17:00'
$data_ws1[17,4] = 505
$data_ws1[17,5] = 60
$data_ws1[17,6] = 'https://show.bilibili.com/platform/detail.html?id=81566'
$data_ws1[17,7] = '//i2.hdslb.com/bfs/openplatform/202401/c4bBhKzu1706685824726.jpeg'
$data_ws1[18,0] = '2024.05.10'
$data_ws1[18,1] = '广州·国际潮宠展—潮流创新宠物展会'
$data_ws1[18,2] = '阅江中路18号 广交会展馆C区'
$data_ws1[18,3] = '2024.05.10 09:30-05.12 18:30'
$data_ws1[18,4] = 8
$data_ws1[18,5] = 36
$data_ws1[18,6] = 'https://show.bilibili.com/platform/detail.html?id=82038'
$data_ws1[18,7] = '//i2.hdslb.com/bfs/openplatform/202402/om8irfxN1708678341525.jpeg'
$data_ws1[19,0] = '2024.05.18'
$data_ws1[19,1] = '广州·恋与深空only'
$data_ws1[19,2] = '大石街石北工业大道644号 巨大创意产业园'
$data_ws1[19,3] = '2024.05.18 10:00-05.18 17:00'
$data_ws1[19,4] = 446
$data_ws1[19,5] = 60
$data_ws1[19,6] = 'https://show.bilibili.com/platform/detail.html?id=81962'
$data_ws1[19,7] = '//i0.hdslb.com/bfs/openplatform/202402/a7aqaXrK1708485268977.jpeg'
$ws1.Range("B2:I21").Value = $data_ws1
# Remove now-stale trailing rows (22-23) so dimension becomes A1:I21
$ws1.Rows("22:23").Delete() | Out-Null

# ---------------- Sheet: 演出 (Performances) ----------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 64
$ws2.Cells.Item(4, 6).Value = 303
$ws2.Cells.Item(6, 6).Value = 13
$ws2.Cells.Item(8, 6).Value = 169
$ws2.Cells.Item(9, 6).Value = 182

# ---------------- Sheet: 本地生活 (Local Life) ----------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 360

# ---------------- Sheet: 全部类型 (All Types) ----------------
$ws4 = $wb.Worksheets.Item("全部类型")
$data_ws4 = New-Object 'object[,]' 32,8
$data_ws4[0,0] = '2024.02.01'
$data_ws4[0,1] = '广州·次元波板糖×线条小狗MALTESE 主题快闪店'
$data_ws4[0,2] = '天河路299号B1层 天河时尚街'
$data_ws4[0,3] = '2024.02.01 00:00-03.01 23:59'
$data_ws4[0,4] = 360
$data_ws4[0,5] = 30
$data_ws4[0,6] = 'https://show.bilibili.com/platform/detail.html?id=81339'
$data_ws4[0,7] = '//i0.hdslb.com/bfs/openplatform/202401/Qbpful951706080847394.png'
$data_ws4[1,0] = '2024.03.02'
$data_ws4[1,1] = '广州·COS STAR次元之夜ACG主题派对'
$data_ws4[1,2] = '南洲路158号2F SD Livehouse'
$data_ws4[1,3] = '2024.03.02 19:00-03.02 22:00'
$data_ws4[1,4] = 64
$data_ws4[1,5] = 118
$data_ws4[1,6] = 'https://show.bilibili.com/platform/detail.html?id=81755'
$data_ws4[1,7] = '//i2.hdslb.com/bfs/openplatform/202402/7PANdxoY1707105412800.png'
$data_ws4[2,0] = '2024.03.02'
$data_ws4[2,1] = '广州·明日方舟ONLY'
$data_ws4[2,2] = '清河东路288号 科尔海悦酒店'
$data_ws4[2,3] = '2024.03.02 10:00-03.02 17:00'
$data_ws4[2,4] = 410
$data_ws4[2,5] = 60
$data_ws4[2,6] = 'https://show.bilibili.com/platform/detail.html?id=80271'
$data_ws4[2,7] = '//i0.hdslb.com/bfs/openplatform/202312/O9z5j7RB1703733578857.jpeg'
$data_ws4[3,0] = '2024.03.03'
$data_ws4[3,1] = '广州·《奥特传奇之希望之光》圆谷正版授权奥特曼系列舞台剧'
$data_ws4[3,2] = '人民北路696号 广州友谊剧院'
$data_ws4[3,3] = '2024.03.03 10:30-03.03 15:10'
$data_ws4[3,4] = 52
$data_ws4[3,5] = 78
$data_ws4[3,6] = 'https://show.bilibili.com/platform/detail.html?id=81023'
$data_ws4[3,7] = '//i0.hdslb.com/bfs/openplatform/202401/r6OKUMAF1705997977504.jpeg'
$data_ws4[4,0] = '2024.03.09'
$data_ws4[4,1] = '广州·HANAPOKO 2024 LIVE'
$data_ws4[4,2] = '海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse'
$data_ws4[4,3] = '2024.03.09 14:00-03.09 15:30'
$data_ws4[4,4] = 303
$data_ws4[4,5] = 380
$data_ws4[4,6] = 'https://show.bilibili.com/platform/detail.html?id=81279'
$data_ws4[4,7] = '//i2.hdslb.com/bfs/openplatform/202401/tMZ1Jp2G1705992352054.jpeg'
$data_ws4[5,0] = '2024.03.09'
$data_ws4[5,1] = '广州·VOCALOID术力口only'
$data_ws4[5,2] = '黄边三横路一街1号 设计殿堂'
$data_ws4[5,3] = '2024.03.09 09:30-03.09 18:00'
$data_ws4[5,4] = 553
$data_ws4[5,5] = 68.8
$data_ws4[5,6] = 'https://show.bilibili.com/platform/detail.html?id=81398'
$data_ws4[5,7] = '//i2.hdslb.com/bfs/openplatform/202401/XpsHJTsC1706160000879.png'
$data_ws4[6,0] = '2024.03.09'
$data_ws4[6,1] = '广州·排球少年.only'
$data_ws4[6,2] = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$data_ws4[6,3] = '2024.03.09 10:00-03.09 17:00'
$data_ws4[6,4] = 868
$data_ws4[6,5] = 60
$data_ws4[6,6] = 'https://show.bilibili.com/platform/detail.html?id=80453'
$data_ws4[6,7] = '//i1.hdslb.com/bfs/openplatform/202401/JYZdnQHl1704341212206.jpeg'
$data_ws4[7,0] = '2024.03.09'
$data_ws4[7,1] = '广州·樱漫动漫嘉年华9.0'
$data_ws4[7,2] = '奥体南路12号 优托邦(奥体旗舰店)'
$data_ws4[7,3] = '2024.03.09 10:00-03.09 17:00'
$data_ws4[7,4] = 595
$data_ws4[7,5] = 58
$data_ws4[7,6] = 'https://show.bilibili.com/platform/detail.html?id=81785'
$data_ws4[7,7] = '//i0.hdslb.com/bfs/openplatform/202402/gBrd7lLX1707203945413.jpeg'
$data_ws4[8,0] = '2024.03.10'
$data_ws4[8,1] = '广州·进击的巨人only'
$data_ws4[8,2] = '机场路1399号广州百信广场二期 李宁运动中心'
$data_ws4[8,3] = '2024.03.10 10:00-03.10 17:00'
$data_ws4[8,4] = 783
$data_ws4[8,5] = 60
$data_ws4[8,6] = 'https://show.bilibili.com/platform/detail.html?id=80454'
$data_ws4[8,7] = '//i0.hdslb.com/bfs/openplatform/202401/m8QkMfFl1704347300282.jpeg'
$data_ws4[9,0] = '2024.03.10'
$data_ws4[9,1] = '广州·马娘only2024'
$data_ws4[9,2] = '黄边三横路一街1号 设计殿堂'
$data_ws4[9,3] = '2024.03.10 09:30-03.10 18:00'
$data_ws4[9,4] = 354
$data_ws4[9,5] = 68.8
$data_ws4[9,6] = 'https://show.bilibili.com/platform/detail.html?id=81632'
$data_ws4[9,7] = '//i0.hdslb.com/bfs/openplatform/202402/svWCXpKm1706776489024.png'
$data_ws4[10,0] = '2024.03.16'
$data_ws4[10,1] = '广州·SISP动漫游戏嘉年华'
$data_ws4[10,2] = '西湾路150号 悦汇城'
$data_ws4[10,3] = '2024.03.16 13:00-03.17 19:00'
$data_ws4[10,4] = 541
$data_ws4[10,5] = 48
$data_ws4[10,6] = 'https://show.bilibili.com/platform/detail.html?id=80624'
$data_ws4[10,7] = '//i2.hdslb.com/bfs/openplatform/202401/Z4Q6Fv8B1704770554777.jpeg'
$data_ws4[11,0] = '2024.03.16'
$data_ws4[11,1] = '广州·YU 7th动漫嘉年华'
$data_ws4[11,2] = '珠江西路8号 高德置地夏广场'
$data_ws4[11,3] = '2024.03.16 10:00-03.17 17:00'
$data_ws4[11,4] = 104
$data_ws4[11,5] = 55
$data_ws4[11,6] = 'https://show.bilibili.com/platform/detail.html?id=81627'
$data_ws4[11,7] = '//i1.hdslb.com/bfs/openplatform/202402/lVqoZMVQ1706775042937.jpeg'
$data_ws4[12,0] = '2024.03.16'
$data_ws4[12,1] = '广州·代号鸢only2.0'
$data_ws4[12,2] = '清河东路288号 科尔海悦酒店'
$data_ws4[12,3] = '2024.03.16 10:00-03.16 21:00'
$data_ws4[12,4] = 1109
$data_ws4[12,5] = 39
$data_ws4[12,6] = 'https://show.bilibili.com/platform/detail.html?id=79828'
$data_ws4[12,7] = '//i0.hdslb.com/bfs/openplatform/202312/RVUVc8oy1702549585918.jpeg'
$data_ws4[13,0] = '2024.03.16'
$data_ws4[13,1] = '广州·原神X星穹铁道X绝区零ONLY'
$data_ws4[13,2] = '洛浦街夏滘西环路1号(厦滘地铁站A口步行290米) 厦喾岭南电商园会展中心'
$data_ws4[13,3] = '2024.03.16 10:00-03.16 17:00'
$data_ws4[13,4] = 558
$data_ws4[13,5] = 60
$data_ws4[13,6] = 'https://show.bilibili.com/platform/detail.html?id=80715'
$data_ws4[13,7] = '//i0.hdslb.com/bfs/openplatform/202401/Lt6ZYvA41704878219924.jpeg'
$data_ws4[14,0] = '2024.03.17'
$data_ws4[14,1] = '广州·三月的幻想演唱会2024「飞越蓝色时刻」'
$data_ws4[14,2] = '恩宁路265号三层、四层自编01 MAO Livehouse广州(永庆坊店)'
$data_ws4[14,3] = '2024.03.17 19:00-03.17 20:30'
$data_ws4[14,4] = 94
$data_ws4[14,5] = 380
$data_ws4[14,6] = 'https://show.bilibili.com/platform/detail.html?id=80870'
$data_ws4[14,7] = '//i1.hdslb.com/bfs/openplatform/202401/8WBT7H6W1705376580145.png'
$data_ws4[15,0] = '2024.03.17'
$data_ws4[15,1] = '广州·梁祝 ·黄河经典名曲大型管弦交响音乐会'
$data_ws4[15,2] = '东风中路299号 广州中山纪念堂'
$data_ws4[15,3] = '2024.03.17 19:30-03.17 21:00'
$data_ws4[15,4] = 13
$data_ws4[15,5] = 75
$data_ws4[15,6] = 'https://show.bilibili.com/platform/detail.html?id=81788'
$data_ws4[15,7] = '//i2.hdslb.com/bfs/openplatform/202402/54YX2MVU1707208994883.jpeg'
$data_ws4[16,0] = '2024.03.23'
$data_ws4[16,1] = '广州·BanG Dream ONLY'
$data_ws4[16,2] = '西环路1号 广州岭南会展中心'
$data_ws4[16,3] = '2024.03.23 10:00-03.23 17:00'
$data_ws4[16,4] = 322
$data_ws4[16,5] = 65
$data_ws4[16,6] = 'https://show.bilibili.com/platform/detail.html?id=81754'
$data_ws4[16,7] = '//i1.hdslb.com/bfs/openplatform/202402/3HJiKSeD1707104926306.jpeg'
$data_ws4[17,0] = '2024.03.23'
$data_ws4[17,1] = '广州·排球少年ONLY'
$data_ws4[17,2] = '机场路1399号广州百信广场二期 李宁运动中心'
$data_ws4[17,3] = '2024.03.23 10:00-03.23 17:00'
$data_ws4[17,4] = 445
$data_ws4[17,5] = 60
$data_ws4[17,6] = 'https://show.bilibili.com/platform/detail.html?id=80716'
$data_ws4[17,7] = '//i0.hdslb.com/bfs/openplatform/202401/IFLvYmxx1704879325152.jpeg'
$data_ws4[18,0] = '2024.03.23'
$data_ws4[18,1] = '广州·春卷饭 十周年  2024  专场演出'
$data_ws4[18,2] = '革新路124号太古仓码头54汇5号仓 太空间Livehouse'
$data_ws4[18,3] = '2024.03.23 20:00-03.23 22:00'
$data_ws4[18,4] = 634
$data_ws4[18,5] = '已售罄'
$data_ws4[18,6] = 'https://show.bilibili.com/platform/detail.html?id=81186'
$data_ws4[18,7] = '//i1.hdslb.com/bfs/openplatform/202401/ho9rIMg21705894649801.jpeg'
$data_ws4[19,0] = '2024.03.24'
$data_ws4[19,1] = '广州·妖都恋与制作人ONLY3.0'
$data_ws4[19,2] = '迎宾大道123号 赛仑吉地大酒店'
$data_ws4[19,3] = '2024.03.24 10:00-03.24 17:30'
$data_ws4[19,4] = 135
$data_ws4[19,5] = 68
$data_ws4[19,6] = 'https://show.bilibili.com/platform/detail.html?id=81715'
$data_ws4[19,7] = '//i1.hdslb.com/bfs/openplatform/202402/SOOWI9wL1708675967102.jpeg'
$data_ws4[20,0] = '2024.03.31'
$data_ws4[20,1] = '广州·KANAKO ITO&AYANE 2024 LIVE'
$data_ws4[20,2] = '奥体南路12号优托邦购物中心 疆进酒Omni Space GZ'
$data_ws4[20,3] = '2024.03.31 19:00-03.31 20:30'
$data_ws4[20,4] = 169
$data_ws4[20,5] = 380
$data_ws4[20,6] = 'https://show.bilibili.com/platform/detail.html?id=81422'
$data_ws4[20,7] = '//i0.hdslb.com/bfs/openplatform/202401/4Y4U8tC01706172039039.jpeg'
$data_ws4[21,0] = '2024.04.06'
$data_ws4[21,1] = '广州·运动番only'
$data_ws4[21,2] = '机场路1399号广州百信广场二期 李宁运动中心'
$data_ws4[21,3] = '2024.04.06 10:00-04.06 17:00'
$data_ws4[21,4] = 287
$data_ws4[21,5] = 55
$data_ws4[21,6] = 'https://show.bilibili.com/platform/detail.html?id=81454'
$data_ws4[21,7] = '//i2.hdslb.com/bfs/openplatform/202401/TBZfwnB41706255329549.jpeg'
$data_ws4[22,0] = '2024.04.13'
$data_ws4[22,1] = '广州·Veni Vidi Vici动漫游戏嘉年华'
$data_ws4[22,2] = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$data_ws4[22,3] = '2024.04.13 10:00-04.13 17:00'
$data_ws4[22,4] = 41
$data_ws4[22,5] = 58
$data_ws4[22,6] = 'https://show.bilibili.com/platform/detail.html?id=81575'
$data_ws4[22,7] = '//i1.hdslb.com/bfs/openplatform/202401/7ir7DZHt1706697841803.jpeg'
$data_ws4[23,0] = '2024.04.13'
$data_ws4[23,1] = '广州·潮娃展WWS'
$data_ws4[23,2] = '西环路1号 广州岭南会展中心'
$data_ws4[23,3] = '2024.04.13 10:00-04.13 17:00'
$data_ws4[23,4] = 64
$data_ws4[23,5] = 48
$data_ws4[23,6] = 'https://show.bilibili.com/platform/detail.html?id=81745'
$data_ws4[23,7] = '//i2.hdslb.com/bfs/openplatform/202402/1SWNaBsA1707100228293.jpeg'
$data_ws4[24,0] = '2024.04.14'
$data_ws4[24,1] = '广州·铃木木乃美 2024 演唱会'
$data_ws4[24,2] = '奥体南路12号优托邦购物中心 疆进酒Omni Space GZ'
$data_ws4[24,3] = '2024.04.14 19:00-04.14 20:30'
$data_ws4[24,4] = 182
$data_ws4[24,5] = 380
$data_ws4[24,6] = 'https://show.bilibili.com/platform/detail.html?id=81911'
$data_ws4[24,7] = '//i0.hdslb.com/bfs/openplatform/202402/rGjpHpAV1708328728461.jpeg'
$data_ws4[25,0] = '2024.04.19'
$data_ws4[25,1] = '广州·动漫钢琴鬼才Kyle Xian互动演奏会'
$data_ws4[25,2] = '人民北路696号 广州友谊剧院'
$data_ws4[25,3] = '2024.04.19 19:30-04.19 21:00'
$data_ws4[25,4] = 39
$data_ws4[25,5] = 64
$data_ws4[25,6] = 'https://show.bilibili.com/platform/detail.html?id=81434'
$data_ws4[25,7] = '//i0.hdslb.com/bfs/openplatform/202401/DJpXVLjd1706236823839.png'
$data_ws4[26,0] = '2024.04.20'
$data_ws4[26,1] = '广州·Arknights Only·夜航星（明日方舟Only)'
$data_ws4[26,2] = '同泰路颐和山庄 颐和大酒店'
$data_ws4[26,3] = '2024.04.20 10:00-04.20 17:00'
$data_ws4[26,4] = 516
$data_ws4[26,5] = 69
$data_ws4[26,6] = 'https://show.bilibili.com/platform/detail.html?id=80282'
$data_ws4[26,7] = '//i2.hdslb.com/bfs/openplatform/202312/gaEHIE1F1703745559785.jpeg'
$data_ws4[27,0] = '2024.04.24'
$data_ws4[27,1] = '广州·今泉爱夏  巡演'
$data_ws4[27,2] = '革新路124号太古仓码头54汇5号仓 太空间Livehouse'
$data_ws4[27,3] = '2024.04.24 20:00-04.24 21:30'
$data_ws4[27,4] = 11
$data_ws4[27,5] = 288
$data_ws4[27,6] = 'https://show.bilibili.com/platform/detail.html?id=81890'
$data_ws4[27,7] = '//i1.hdslb.com/bfs/openplatform/202402/YJENeaUi1708313389899.jpeg'
$data_ws4[28,0] = '2024.04.28'
$data_ws4[28,1] = ' 广州·夏川里美 2024 巡回演唱会 出道 25 周年纪念专场'
$data_ws4[28,2] = '中山纪念堂 中山纪念堂'
$data_ws4[28,3] = '2024.04.28 19:30-04.28 21:30'
$data_ws4[28,4] = 17
$data_ws4[28,5] = 280
$data_ws4[28,6] = 'https://show.bilibili.com/platform/detail.html?id=81068'
$data_ws4[28,7] = '//i0.hdslb.com/bfs/openplatform/202401/pXznRv8G1705633441713.jpeg'
$data_ws4[29,0] = '2024.05.05'
$data_ws4[29,1] = '广州·第八届萌物语动漫嘉年华'
$data_ws4[29,2] = '洛浦街厦滘西环路1号 岭南会展中心'
$data_ws4[29,3] = '2024.05.05 10:00-05.05 17:00'
$data_ws4[29,4] = 505
$data_ws4[29,5] = 60
$data_ws4[29,6] = 'https://show.bilibili.com/platform/detail.html?id=81566'
$data_ws4[29,7] = '//i2.hdslb.com/bfs/openplatform/202401/c4bBhKzu1706685824726.jpeg'
$data_ws4[30,0] = '2024.05.10'
$data_ws4[30,1] = '广州·国际潮宠展—潮流创新宠物展会'
$data_ws4[30,2] = '阅江中路18号 广交会展馆C区'
$data_ws4[30,3] = '2024.05.10 09:30-05.12 18:30'
$data_ws4[30,4] = 8
$data_ws4[30,5] = 36
$data_ws4[30,6] = 'https://show.bilibili.com/platform/detail.html?id=82038'
$data_ws4[30,7] = '//i2.hdslb.com/bfs/openplatform/202402/om8irfxN1708678341525.jpeg'
$data_ws4[31,0] = '2024.05.18'
$data_ws4[31,1] = '广州·恋与深空only'
$data_ws4[31,2] = '大石街石北工业大道644号 巨大创意产业园'
$data_ws4[31,3] = '2024.05.18 10:00-05.18 17:00'
$data_ws4[31,4] = 446
$data_ws4[31,5] = 60
$data_ws4[31,6] = 'https://show.bilibili.com/platform/detail.html?id=81962'
$data_ws4[31,7] = '//i0.hdslb.com/bfs/openplatform/202402/a7aqaXrK1708485268977.jpeg'
$ws4.Range("B2:I33").Value = $data_ws4
# Remove now-stale trailing rows (34-35) so dimension becomes A1:I33
$ws4.Rows("34:35") .Delete() | Out-Null

Write-Host "Edit applied."
